$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Lucene" (sheet1) - fill in previously-missing Binary Relevance
# results for Config 6 (Count Vectorizer + TFIDF + ngram(2) + POS / Random
# Forest row, and the whole Count Vectorizer + TFIDF + ngram(1) + POS block).
# ---------------------------------------------------------------------------
$wsLucene = $wb.Worksheets.Item("Lucene")

# Row 55: Random Forest result row for the "ngram(2) + POS" configuration
$wsLucene.Range("C55").Value = "0.170 0.728 0.543 0.263 0.664"
$wsLucene.Range("D55").Value = "0.889 0.730 0.857 0.896 0.936"
$wsLucene.Range("E55").Value = "0.093 0.620 0.376 0.151 0.498"
$wsLucene.Range("F55").Value = "0.903 0.794 0.848 0.899 0.955"

# Row 57: Logistic Regression result row for "ngram(1) + POS" configuration
$wsLucene.Range("C57").Value = "0.457 0.666 0.494 0.480 0.720 "
$wsLucene.Range("D57").Value = "0.118 0.672 0.665 0.137 0.848"
$wsLucene.Range("E57").Value = "0.340 0.541 0.333 0.366 0.565 "
$wsLucene.Range("F57").Value = "0.661 0.755 0.815 0.659 0.955"

# Row 58: Multinomial Naive Bayes result row for "ngram(1) + POS" configuration
$wsLucene.Range("C58").Value = "0.495 0.728 0.750 0.664 0.835"
$wsLucene.Range("D58").Value = "0.319 0.549 0.444 0.249 0.546"
$wsLucene.Range("E58").Value = "0.340 0.816 0.783 0.585 0.749"
$wsLucene.Range("F58").Value = "0.853 0.711 0.734 0.747 0.926"

# Row 59: Support Vector Machines result row for "ngram(1) + POS" configuration
$wsLucene.Range("C59").Value = "0.424 0.718 0.699 0.519 0.779"
$wsLucene.Range("D59").Value = "0.336 0.668 0.633 0.420 0.732 "
$wsLucene.Range("E59").Value = "0.274 0.626 0.569 0.359 0.647 "
$wsLucene.Range("F59").Value = "0.866 0.768 0.831 0.868 0.950"

# Row 60: Decision Tree result row for "ngram(1) + POS" configuration
$wsLucene.Range("C60").Value = "0.053 0.606 0.486 0.156 0.635"
$wsLucene.Range("D60").Value = "0.636 0.599 0.668 0.453 0.735"
$wsLucene.Range("E60").Value = "0.027 0.475 0.326 0.085 0.469"
$wsLucene.Range("F60").Value = "0.895 0.715 0.814 0.882 0.941"

# Row 61: Random Forest result row for "ngram(1) + POS" configuration
$wsLucene.Range("C61").Value = "0.163 0.734 0.520 0.268 0.677 "
$wsLucene.Range("D61").Value = "1.000 0.752 0.865 0.917 0.955"
$wsLucene.Range("E61").Value = "0.089 0.622 0.354 0.155 0.512 "
$wsLucene.Range("F61").Value = "0.904 0.803 0.844 0.900 0.957"

# Update the selection/active cell on the Lucene sheet and deselect its tab.
$wsLucene.Range("C63").Select()

# ---------------------------------------------------------------------------
# Sheet view / active tab bookkeeping: the workbook's active sheet moves from
# "Lucene" to "Thunderbird", and the Thunderbird sheet's remembered selection
# moves from A45 to B49.
# ---------------------------------------------------------------------------
$wsThunderbird = $wb.Worksheets.Item("Thunderbird")
$wsThunderbird.Range("B49").Select()
$wsThunderbird.Activate()

Write-Host "done"
